$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# Update the Tags value for the first equipment row (row 2) to be more
# descriptive, adding a new shared string "Medical, ICU".
$ws.Range("C2").Value = "Medical, ICU"

# Move/update the active selection to the edited cell.
$ws.Range("C2").Select()
